$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 120; existing rows 120-139 shift down to 121-140.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with its data (weekly series entry).
$ws.Cells.Item(120, 1).Value = 3
$ws.Cells.Item(120, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44505
$ws.Cells.Item(120, 5).Value = 5
$ws.Cells.Item(120, 6).Value = 100112010
$ws.Cells.Item(120, 7).Value = "Achicoria"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 115
$ws.Cells.Item(120, 11).Value = 6000
$ws.Cells.Item(120, 12).Value = 6500
$ws.Cells.Item(120, 13).Value = 6239
$ws.Cells.Item(120, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(120, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(120, 16).Value = 390
$ws.Cells.Item(120, 17).Value = 16
$ws.Cells.Item(120, 18).Value = "Hortaliza"
